# Remove Stop codons (TAG, TAA, TGA) from the "preferentially_used" (column C)
# and "non_preferentially_used" (column D) codon groups on both worksheets,
# shifting the remaining values up and clearing the now-empty trailing cells.

$wb = $excel.ActiveWorkbook

foreach ($ws in @($wb.Worksheets.Item("Human"), $wb.Worksheets.Item("Mosquito"))) {

    # --- Column C ("preferentially_used"): remove TAG at row 12, shift rows 13-15 up,
    #     clear the newly-vacated row 15 ---
    $colC = @()
    for ($r = 12; $r -le 15; $r++) {
        $colC += ,$ws.Range("C$r").Value2
    }
    # Drop the stop codon (TAG) which is always the first of this block
    $colC = $colC[1..($colC.Length - 1)]

    for ($i = 0; $i -lt $colC.Length; $i++) {
        $ws.Range("C" + (12 + $i)).Value = $colC[$i]
    }
    $ws.Range("C15").Value = $null

    # --- Column D ("non_preferentially_used"): remove TAA at row 17 and TGA at row 19,
    #     shift row 18 up into row 17, then delete the now-empty rows 18 and 19 ---
    $ws.Range("D17").Value = $ws.Range("D18").Value2
    $ws.Range("A18:F19").Delete()
}
